$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marksheet "Corr/total marks" update (row 11 = "Marking", row 12 = "Total")
# Right-answer count for the "Marking" row: 3 -> 5
$ws.Range("B11").Value = 5

# Right-answer count for the "Total" row: 42 -> 70
$ws.Range("B12").Value = 70

# Correct/Total marks summary string for the "Total" row: "28/84" -> "70/140"
$ws.Range("E12").Value = "70/140"
